$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.903.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.420.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.11"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.404.14"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.05%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.844.66"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.765.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.392.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.30%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.68"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.93"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.06"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.88"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "584.76"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.43"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -9.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.524.54"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.93%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.53%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.21%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.02"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.15"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.15"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.41%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0286"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +15.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.19"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.589"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.55"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0504"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.12%  "
